$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: remove the whole paragraph "IMPORT color output library
# (colorama)" (merges away, leaving the following blank-space paragraph
# intact).
# ---------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*IMPORT color output library (colorama)*") {
        $p.Range.Delete()
        $found = $true
        break
    }
}
Write-Host "Removed colorama paragraph: $found"

# ---------------------------------------------------------------------
# Change 2: split "    ELSE IF guess is in word_letters:" so that
# "word_letters" (and the trailing ":") become their own runs.
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("word_letters:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.MoveEnd(1, -1)
    $r.Bold = 1
    $r.Bold = 0
}
Write-Host "Split word_letters run: $($r.Find.Found)"

# ---------------------------------------------------------------------
# Change 3: "        PRINT incorrect guess message in red" ->
# "        PRINT incorrect guess message"
# ---------------------------------------------------------------------
$replaced = $d.Content.Find.Execute("PRINT incorrect guess message in red", $true, $false, $false, $false, $false, $true, 1, $false, "PRINT incorrect guess message", 2)
Write-Host "Replaced incorrect guess text: $replaced"

$full = $d.Content.Text
$needle = "        PRINT incorrect guess message"
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $rng = $d.Range($idx, $idx + $needle.Length)
    $rng.Bold = 1
    $rng.Bold = 0
}

# ---------------------------------------------------------------------
# Change 4: "IF lives > 0: PRINT congratulatory message in green
# INCREMENT score[...]..." -> split into
#   "IF lives > 0: PRINT congratulatory message"
#   ","
#   " INCREMENT score[...] ..."
# ---------------------------------------------------------------------
$replaced2 = $d.Content.Find.Execute("congratulatory message in green INCREMENT", $true, $false, $false, $false, $false, $true, 1, $false, "congratulatory message, INCREMENT", 2)
Write-Host "Replaced congratulatory text: $replaced2"

$full = $d.Content.Text
$prefix = "IF lives > 0: PRINT congratulatory message,"
$idx = $full.IndexOf($prefix)
if ($idx -ge 0) {
    $wholeText = "IF lives > 0: PRINT congratulatory message, INCREMENT score[`"wins`"] ELSE: PRINT game over message and reveal word INCREMENT score[`"losses`"]"
    $wholeRng = $d.Range($idx, $idx + $wholeText.Length)
    $wholeRng.Bold = 1
    $wholeRng.Bold = 0

    $commaIdx = $idx + "IF lives > 0: PRINT congratulatory message".Length
    $commaRng = $d.Range($commaIdx, $commaIdx + 1)
    Write-Host "Comma char: [$($commaRng.Text)]"
    $commaRng.Bold = 1
    $commaRng.Bold = 0
}

Write-Host "Done"
